$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric; keep them as text
# to match the source data (all originally stored as inline strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.819.12'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '1.631.52'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '215.59'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '0.5058'
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("D9").Value = '0.06436'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '19.49'
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").Value = '0.07785'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '4.274'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.633.37'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.856.00'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '0.5605'
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").Value = '0.0₅7585'
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").Value = '63.06'
$ws.Range("E17").Value = '  -2.00%  '
$ws.Range("D18").Value = '25.831.49'
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '194.71'
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("D22").Value = '9.852'
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = '6.028'
$ws.Range("E23").Value = '  -2.21%  '
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '1.798'
$ws.Range("E25").Value = '  -5.24%  '
$ws.Range("D26").Value = '141.15'
$ws.Range("D27").Value = '0.1278'
$ws.Range("E27").Value = '  +1.30%  '
$ws.Range("D28").Value = '6.752'
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("D29").Value = '15.43'
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").Value = '1.239'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").Value = '0.04881'
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").Value = '3.222'
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").Value = '1.552'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").Value = '2.378'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '0.8967'
$ws.Range("E36").Value = '  -2.56%  '
$ws.Range("D37").Value = '2.569'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '1.129.00'
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").Value = '0.5498'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").Value = '0.01562'
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("D41").Value = '0.9936'
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("D42").Value = '5.536'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("D43").Value = '0.7977'
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").Value = '97.33'
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").Value = '1.781.79'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  -3.77%  '
$ws.Range("D47").Value = '0.4441'
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").Value = '55.33'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = '0.05055'
$ws.Range("D50").Value = '7.659'
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").Value = '0.9978'
$ws.Range("E51").Value = '  -0.41%  '
